$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 48: id 47 - The River Merchant NPC
$ws.Cells.Item(48, 1).Value = 47
$ws.Cells.Item(48, 2).Value = "TheRiverMerchant"
$ws.Cells.Item(48, 3).Value = "The River Merchant"
$ws.Cells.Item(48, 4).Value = 2
$ws.Cells.Item(48, 5).Value = "Twisted Memories"
$ws.Cells.Item(48, 9).Value = 1280
$ws.Cells.Item(48, 10).Value = 272

# Row 49: id 48 - Child of Alchemy NPC
$ws.Cells.Item(49, 1).Value = 48
$ws.Cells.Item(49, 2).Value = "ChildofAlchemy"
$ws.Cells.Item(49, 3).Value = "Child of Alchemy"
$ws.Cells.Item(49, 4).Value = 2
$ws.Cells.Item(49, 5).Value = "Twisted Memories"
$ws.Cells.Item(49, 9).Value = 1408
$ws.Cells.Item(49, 10).Value = 64
